# Add a new "intervention_type" column (K) to the clinical trials list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell K1, with same style (bold / centered / bordered) as the other headers.
$ws.Range("K1").Value = "intervention_type"
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null # xlPasteFormats

# Values for K2:K26, in row order.
$interventionTypes = @(
    "DIETARY_SUPPLEMENT", # row 2
    "DIETARY_SUPPLEMENT", # row 3
    "DIETARY_SUPPLEMENT", # row 4
    "DIETARY_SUPPLEMENT", # row 5
    "DIETARY_SUPPLEMENT", # row 6
    "DIETARY_SUPPLEMENT", # row 7
    "DIETARY_SUPPLEMENT", # row 8
    "DEVICE",             # row 9
    "OTHER",              # row 10
    "DEVICE",             # row 11
    "DEVICE",             # row 12
    "OTHER",              # row 13
    "DIETARY_SUPPLEMENT", # row 14
    "OTHER",              # row 15
    "BEHAVIORAL",         # row 16
    "OTHER",              # row 17
    "OTHER",              # row 18
    "OTHER",              # row 19
    "DIETARY_SUPPLEMENT", # row 20
    "OTHER",              # row 21
    "DIETARY_SUPPLEMENT", # row 22
    "DRUG",               # row 23
    "DIETARY_SUPPLEMENT", # row 24
    "DIETARY_SUPPLEMENT", # row 25
    "DIETARY_SUPPLEMENT"  # row 26
)

$rowCount = $interventionTypes.Length
$arr = New-Object 'object[,]' $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $arr[$i, 0] = $interventionTypes[$i]
}

$ws.Range("K2:K26").Value = $arr

Write-Output "Done: set K1:K26 intervention_type column"
